$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26: becomes the "Violgubbe" (Gomphus clavatus) record
$ws.Range("A26").Value = 112093186
$ws.Range("B26").Value = 89043
$ws.Range("D26").Value = "VU"
$ws.Range("E26").Value = 720
$ws.Range("F26").Value = "Violgubbe"
$ws.Range("G26").Value = "Gomphus clavatus"
$ws.Range("H26").Value = "(Pers.) Gray"
$ws.Range("Q26").Value = 639205
$ws.Range("R26").Value = 6701016
$ws.Range("AC26").Value = "Till stor del barkborredödat bestånd intill stort kalhygge"
$ws.Range("AD26").Value = $true

# Row 27: only the Taxonsorteringsordning (B) value changes
$ws.Range("B27").Value = 90821

# Row 28: becomes the "Strimspindling" (Cortinarius glaucopus) record
$ws.Range("A28").Value = 112093190
$ws.Range("B28").Value = 85331
$ws.Range("D28").Value = "LC"
$ws.Range("E28").Value = 3624
$ws.Range("F28").Value = "Strimspindling"
$ws.Range("G28").Value = "Cortinarius glaucopus"
$ws.Range("H28").Value = "(Schaeff. : Fr.) Fr."
$ws.Range("P28").Value = "Fagerdal, Upl"
$ws.Range("Q28").Value = 639180
$ws.Range("R28").Value = 6701165
$ws.Range("AC28").Value = ""
$ws.Range("AD28").Value = $false

# Row 29: only the Taxonsorteringsordning (B) value changes
$ws.Range("B29").Value = 89317
